$d = $word.ActiveDocument

# 1. Title: merge " " + "ARCHive" + " Format Risk Report" into one run's text
#    (removes the proofErr spell-check markers around "ARCHive" along the way)
$d.Content.Find.Execute(" ARCHive Format Risk Report", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " ARCHive Format Risk Report", 2) | Out-Null

# 2. Overview paragraph: merge "ARCHive" + " as of 2023-1" into one run's text
$d.Content.Find.Execute("ARCHive as of 2023-1", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "ARCHive as of 2023-1", 2) | Out-Null

# 3. Risk Profile paragraph: insert "NARA " before "risk level for the entire department."
$d.Content.Find.Execute(" at each risk level for the entire department.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " at each NARA risk level for the entire department.", 2) | Out-Null

# 4. Risk Change paragraph: insert "NARA " before "risk level from the previous analysis and "
$d.Content.Find.Execute("of formats at each risk level from the previous analysis and ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "of formats at each NARA risk level from the previous analysis and ", 2) | Out-Null

# 5. Risk By Collection paragraph: insert "NARA " before "risk level in each collection" and
#    also collapse the proofErr-wrapped "Collection_Risk_Levels" run boundary into a single run.
$d.Content.Find.Execute("The percentage of formats at each risk level in each collection. The goal is to visualize how many collections have a high percentage of risk vs. low percentage of risk. The risk report spreadsheet (" + [char]8220 + "Collection_Risk_Levels" + [char]8221 + ") has the risk data for individual collections.", `
                         $false, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "The percentage of formats at each NARA risk level in each collection. The goal is to visualize how many collections have a high percentage of risk vs. low percentage of risk. The risk report spreadsheet (" + [char]8220 + "Collection_Risk_Levels" + [char]8221 + ") has the risk data for individual collections.", 2) | Out-Null
